$d = $word.ActiveDocument

# Locate the target paragraph (1-based index) containing the original
# sentence that is being replaced/expanded.
$found = 0
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Sau đây là kết quả khi chạy các file ở trên*") {
        $found = $idx
        break
    }
}
if ($found -eq 0) {
    throw "Target paragraph not found"
}

$sentencePara = $d.Paragraphs($found)
$blankPara = $d.Paragraphs($found - 1)

# 1. The blank paragraph right before the sentence gains a first-line
#    indent of 360 twips (18pt) -- matches the <w:ind w:firstLine="360"/>
#    added to its pPr.
$blankPara.Format.FirstLineIndent = 18

# 2. Insert a brand-new paragraph right after it holding the lead-in
#    sentence "Sử dụng lệnh sau: ".
$blankPara.Range.InsertParagraphAfter()
$introPara = $d.Paragraphs($found)
$introPara.Range.Text = "Sử dụng lệnh sau: "
$introPara.Range.Font.Name = "Times New Roman"
$introPara.Range.Font.NameBi = "Times New Roman"
$introPara.Format.FirstLineIndent = 18

# 3. Insert another new, centered paragraph after that one holding the
#    shell command itself.
$introPara.Range.InsertParagraphAfter()
$cmdPara = $d.Paragraphs($found + 1)
$cmdPara.Range.Text = "python3 /opt/run_crud/py"
$cmdPara.Range.Font.Name = "Times New Roman"
$cmdPara.Range.Font.NameBi = "Times New Roman"
$cmdPara.Format.FirstLineIndent = 18
$cmdPara.Format.Alignment = 1

# 4. The original sentence paragraph is rewritten as two runs:
#    "để chạy file run_crud. " followed by
#    "Sau đây là kết quả khi chạy file:". To get genuinely distinct
#    <w:r> runs (rather than one coalesced run) build the two pieces as
#    separate paragraphs first -- each paragraph-level Range.Font.Name/
#    NameBi assignment reliably stamps the full rFonts (ascii/hAnsi/cs)
#    -- then join them back into a single paragraph by deleting the
#    paragraph mark between them.
$sentencePara = $d.Paragraphs($found + 2)
$sentencePara.Range.Text = "để chạy file run_crud. "
$sentencePara.Range.Font.Name = "Times New Roman"
$sentencePara.Range.Font.NameBi = "Times New Roman"
$sentencePara.Format.FirstLineIndent = 18
$sentencePara.Format.Alignment = 0

$sentencePara.Range.InsertParagraphAfter()
$tailPara = $d.Paragraphs($found + 3)
$tailPara.Range.Text = "Sau đây là kết quả khi chạy file:"
$tailPara.Range.Font.Name = "Times New Roman"
$tailPara.Range.Font.NameBi = "Times New Roman"
$tailPara.Format.FirstLineIndent = 18

$joinMark = $d.Range($sentencePara.Range.End - 1, $sentencePara.Range.End)
$joinMark.Delete()
